$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("PIA","CONCURSO APOYO CENTROS DE EXCELENCIA CON FINANCIAMIENTO BASAL 2022","AFB220001","2022","2022","DORA ALTBIR DRULLINSKY","MUJER","INSTITUCION ASOCIADA","UNIVERSIDAD DE CHILE","60910000-1","CHILE","REGION METROPOLITANA DE SANTIAGO"),
    @("PIA","CONCURSO APOYO CENTROS DE EXCELENCIA CON FINANCIAMIENTO BASAL 2022","AFB220001","2022","2022","DORA ALTBIR DRULLINSKY","MUJER","INSTITUCION ALBERGANTE","UNIVERSIDAD DE SANTIAGO DE CHILE","60911000-7","CHILE","NO REGISTRA REGION"),
    @("PIA","CONCURSO APOYO CENTROS DE EXCELENCIA CON FINANCIAMIENTO BASAL 2022","AFB220001","2022","2022","DORA ALTBIR DRULLINSKY","MUJER","INSTITUCION ASOCIADA","UNIVERSIDAD TECNICA FEDERICO SANTA MARIA","81668700-4","CHILE","REGION DE VALPARAISO"),
    @("PIA","CONCURSO APOYO CENTROS DE EXCELENCIA CON FINANCIAMIENTO BASAL 2022","AFB220002","2022","2022","JAVIER RUIZ DEL SOLAR SAN MARTIN","HOMBRE","INSTITUCION ALBERGANTE","UNIVERSIDAD DE CHILE","60910000-1","CHILE","REGION METROPOLITANA DE SANTIAGO"),
    @("PIA","CONCURSO APOYO CENTROS DE EXCELENCIA CON FINANCIAMIENTO BASAL 2022","AFB220003","2022","2022","LEONARDO BASSO SOTZ","HOMBRE","INSTITUCION ALBERGANTE","INSTITUTO SISTEMAS COMPLEJOS DE INGENIERIA","65908140-7","CHILE","REGION METROPOLITANA DE SANTIAGO"),
    @("PIA","CONCURSO APOYO CENTROS DE EXCELENCIA CON FINANCIAMIENTO BASAL 2022","AFB220004","2022","2022","IVAN SCHMIDT ANDRADE","HOMBRE","INSTITUCION ALBERGANTE","UNIVERSIDAD TECNICA FEDERICO SANTA MARIA","81668700-4","CHILE","REGION DE VALPARAISO"),
    @("FONDAP","CONCURSO DE APOYO A CENTROS DE EXCELENCIA EN INVESTIGACION EN AREAS PRIORITARIAS - FONDAP 2022","1522A0001","2022","2022","RENE GARREAUD SALAZAR","HOMBRE","INSTITUCION PATROCINANTE","UNIVERSIDAD DE CHILE","60910000-1","CHILE","REGION METROPOLITANA DE SANTIAGO"),
    @("FONDAP","CONCURSO DE APOYO A CENTROS DE EXCELENCIA EN INVESTIGACION EN AREAS PRIORITARIAS - FONDAP 2022","1522A0001","2022","2022","RENE GARREAUD SALAZAR","HOMBRE","INSTITUCION ASOCIADA","UNIVERSIDAD AUSTRAL DE CHILE","81380500-6","CHILE","REGION DE LOS RIOS"),
    @("FONDAP","CONCURSO DE APOYO A CENTROS DE EXCELENCIA EN INVESTIGACION EN AREAS PRIORITARIAS - FONDAP 2022","1522A0001","2022","2022","RENE GARREAUD SALAZAR","HOMBRE","INSTITUCION ASOCIADA","UNIVERSIDAD DE CONCEPCION","81494400-K","CHILE","REGION DEL BIO-BIO"),
    @("FONDAP","CONCURSO DE APOYO A CENTROS DE EXCELENCIA EN INVESTIGACION EN AREAS PRIORITARIAS - FONDAP 2022","1522A0002","2022","2022","WALDO BUSTAMANTE GOMEZ","HOMBRE","INSTITUCION ASOCIADA","UNIVERSIDAD DE CONCEPCION","81494400-K","CHILE","REGION DEL BIO-BIO"),
    @("FONDAP","CONCURSO DE APOYO A CENTROS DE EXCELENCIA EN INVESTIGACION EN AREAS PRIORITARIAS - FONDAP 2022","1522A0002","2022","2022","WALDO BUSTAMANTE GOMEZ","HOMBRE","INSTITUCION PATROCINANTE","PONTIFICIA UNIVERSIDAD CATOLICA DE CHILE","81698900-0","CHILE","REGION METROPOLITANA DE SANTIAGO"),
    @("FONDAP","CONCURSO DE APOYO A CENTROS DE EXCELENCIA EN INVESTIGACION EN AREAS PRIORITARIAS - FONDAP 2022","1522A0003","2022","2022","PEDRO MEGE ROSSO","HOMBRE","INSTITUCION ASOCIADA","UNIVERSIDAD DIEGO PORTALES","70990700-K","CHILE","REGION METROPOLITANA DE SANTIAGO"),
    @("FONDAP","CONCURSO DE APOYO A CENTROS DE EXCELENCIA EN INVESTIGACION EN AREAS PRIORITARIAS - FONDAP 2022","1522A0003","2022","2022","PEDRO MEGE ROSSO","HOMBRE","INSTITUCION ASOCIADA","UNIVERSIDAD ACADEMIA DE HUMANISMO CRISTIANO","71470400-1","CHILE","REGION METROPOLITANA DE SANTIAGO"),
    @("FONDAP","CONCURSO DE APOYO A CENTROS DE EXCELENCIA EN INVESTIGACION EN AREAS PRIORITARIAS - FONDAP 2022","1522A0003","2022","2022","PEDRO MEGE ROSSO","HOMBRE","INSTITUCION PATROCINANTE","PONTIFICIA UNIVERSIDAD CATOLICA DE CHILE","81698900-0","CHILE","REGION METROPOLITANA DE SANTIAGO"),
    @("FONDAP","CONCURSO DE APOYO A CENTROS DE EXCELENCIA EN INVESTIGACION EN AREAS PRIORITARIAS - FONDAP 2022","1522A0004","2022","2022","RENATO QUIÑONES BERGERET","HOMBRE","INSTITUCION ASOCIADA","UNIVERSIDAD ANDRES BELLO","71540100-2","CHILE","REGION METROPOLITANA DE SANTIAGO"),
    @("FONDAP","CONCURSO DE APOYO A CENTROS DE EXCELENCIA EN INVESTIGACION EN AREAS PRIORITARIAS - FONDAP 2022","1522A0004","2022","2022","RENATO QUIÑONES BERGERET","HOMBRE","INSTITUCION ASOCIADA","UNIVERSIDAD AUSTRAL DE CHILE","81380500-6","CHILE","REGION DE LOS RIOS"),
    @("FONDAP","CONCURSO DE APOYO A CENTROS DE EXCELENCIA EN INVESTIGACION EN AREAS PRIORITARIAS - FONDAP 2022","1522A0004","2022","2022","RENATO QUIÑONES BERGERET","HOMBRE","INSTITUCION PATROCINANTE","UNIVERSIDAD DE CONCEPCION","81494400-K","CHILE","REGION DEL BIO-BIO"),
    @("FONDAP","CONCURSO DE APOYO A CENTROS DE EXCELENCIA EN INVESTIGACION EN AREAS PRIORITARIAS - FONDAP 2022","1522A0005","2022","2022","RODRIGO CIENFUEGOS CARRASCO","HOMBRE","INSTITUCION ASOCIADA","UNIVERSIDAD ANDRES BELLO","71540100-2","CHILE","REGION METROPOLITANA DE SANTIAGO"),
    @("FONDAP","CONCURSO DE APOYO A CENTROS DE EXCELENCIA EN INVESTIGACION EN AREAS PRIORITARIAS - FONDAP 2022","1522A0005","2022","2022","RODRIGO CIENFUEGOS CARRASCO","HOMBRE","INSTITUCION ASOCIADA","UNIVERSIDAD CATOLICA DEL NORTE","81518400-9","CHILE","REGION DE ANTOFAGASTA"),
    @("FONDAP","CONCURSO DE APOYO A CENTROS DE EXCELENCIA EN INVESTIGACION EN AREAS PRIORITARIAS - FONDAP 2022","1522A0005","2022","2022","RODRIGO CIENFUEGOS CARRASCO","HOMBRE","INSTITUCION ASOCIADA","UNIVERSIDAD TECNICA FEDERICO SANTA MARIA","81668700-4","CHILE","REGION DE VALPARAISO"),
    @("FONDAP","CONCURSO DE APOYO A CENTROS DE EXCELENCIA EN INVESTIGACION EN AREAS PRIORITARIAS - FONDAP 2022","1522A0005","2022","2022","RODRIGO CIENFUEGOS CARRASCO","HOMBRE","INSTITUCION PATROCINANTE","PONTIFICIA UNIVERSIDAD CATOLICA DE CHILE","81698900-0","CHILE","REGION METROPOLITANA DE SANTIAGO"),
    @("FONDAP","CONCURSO DE APOYO A CENTROS DE EXCELENCIA EN INVESTIGACION EN AREAS PRIORITARIAS - FONDAP 2022","1522A0006","2022","2022","CLAUDIA RAHMANN ZUÑIGA","MUJER","INSTITUCION PATROCINANTE","UNIVERSIDAD DE CHILE","60910000-1","CHILE","REGION METROPOLITANA DE SANTIAGO"),
    @("FONDAP","CONCURSO DE APOYO A CENTROS DE EXCELENCIA EN INVESTIGACION EN AREAS PRIORITARIAS - FONDAP 2022","1522A0006","2022","2022","CLAUDIA RAHMANN ZUÑIGA","MUJER","INSTITUCION ASOCIADA","FUNDACION FRAUNHOFER CHILE RESEARCH","65033192-3","CHILE","REGION METROPOLITANA DE SANTIAGO"),
    @("FONDAP","CONCURSO DE APOYO A CENTROS DE EXCELENCIA EN INVESTIGACION EN AREAS PRIORITARIAS - FONDAP 2022","1522A0006","2022","2022","CLAUDIA RAHMANN ZUÑIGA","MUJER","INSTITUCION ASOCIADA","UNIVERSIDAD DE TARAPACA","70770800-K","CHILE","REGION DE ARICA Y PARINACOTA"),
    @("FONDAP","CONCURSO DE APOYO A CENTROS DE EXCELENCIA EN INVESTIGACION EN AREAS PRIORITARIAS - FONDAP 2022","1522A0006","2022","2022","CLAUDIA RAHMANN ZUÑIGA","MUJER","INSTITUCION ASOCIADA","UNIVERSIDAD DE ANTOFAGASTA","70791800-4","CHILE","REGION DE ANTOFAGASTA"),
    @("FONDAP","CONCURSO DE APOYO A CENTROS DE EXCELENCIA EN INVESTIGACION EN AREAS PRIORITARIAS - FONDAP 2022","1522A0006","2022","2022","CLAUDIA RAHMANN ZUÑIGA","MUJER","INSTITUCION ASOCIADA","UNIVERSIDAD ADOLFO IBANEZ","71543200-5","CHILE","REGION METROPOLITANA DE SANTIAGO"),
    @("FONDAP","CONCURSO DE APOYO A CENTROS DE EXCELENCIA EN INVESTIGACION EN AREAS PRIORITARIAS - FONDAP 2022","1522A0006","2022","2022","CLAUDIA RAHMANN ZUÑIGA","MUJER","INSTITUCION ASOCIADA","UNIVERSIDAD DE CONCEPCION","81494400-K","CHILE","REGION DEL BIO-BIO"),
    @("FONDAP","CONCURSO DE APOYO A CENTROS DE EXCELENCIA EN INVESTIGACION EN AREAS PRIORITARIAS - FONDAP 2022","1522A0006","2022","2022","CLAUDIA RAHMANN ZUÑIGA","MUJER","INSTITUCION ASOCIADA","UNIVERSIDAD TECNICA FEDERICO SANTA MARIA","81668700-4","CHILE","REGION DE VALPARAISO"),
    @("FONDAP","CONCURSO DE APOYO A CENTROS DE EXCELENCIA EN INVESTIGACION EN AREAS PRIORITARIAS - FONDAP 2022","1522A0006","2022","2022","CLAUDIA RAHMANN ZUÑIGA","MUJER","INSTITUCION ASOCIADA","PONTIFICIA UNIVERSIDAD CATOLICA DE CHILE","81698900-0","CHILE","REGION METROPOLITANA DE SANTIAGO"),
    @("FONDAP","APOYO AL FORTALECIMIENTO A LAS CAPACIDADES PARA I+D 2022 MODALIDAD CENTROS TECNOLOGICOS PARA LA INNOVACION Y CENTROS DE EXCELENCIA INTERNACIONAL","CTI220001","2022","2022","MAY CHOMALI GARIB","MUJER","CENTRO","CENTRO NACIONAL DE SISTEMAS DE                    INFORMACION EN SALUD","65165006-2","CHILE","REGION METROPOLITANA DE SANTIAGO"),
    @("FONDAP","APOYO AL FORTALECIMIENTO A LAS CAPACIDADES PARA I+D 2022 MODALIDAD CENTROS TECNOLOGICOS PARA LA INNOVACION Y CENTROS DE EXCELENCIA INTERNACIONAL","CTI220002","2022","2022","NAYAT SANCHEZ PI","MUJER","CENTRO","FUNDACION INRIA CHILE","65058044-3","CHILE","REGION METROPOLITANA DE SANTIAGO")
)

$startRow = 1238
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $row = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = [int]$row[3]
    $ws.Cells.Item($r, 5).Value = [int]$row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $ws.Cells.Item($r, 9).Value = $row[8]
    $ws.Cells.Item($r, 10).Value = $row[9]
    $ws.Cells.Item($r, 11).Value = $row[10]
    $ws.Cells.Item($r, 12).Value = $row[11]
}

# Update the _FilterDatabase defined name to reflect the data range
for ($i = 1; $i -le $wb.Names.Count; $i++) {
    $n = $wb.Names.Item($i)
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "=SCIA_INSTITUCIONES!`$A`$1:`$L`$1237"
    }
}

# Update view selection to mirror final state
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 1245
$ws.Range("A1269").Select()
